$wb = $excel.ActiveWorkbook

# --- SEC_Comm sheet ---
$wsComm = $wb.Sheets.Item("SEC_Comm")

$wsComm.Range("B7").Style = "Normal"
$wsComm.Range("B7").Value = "NRG"
$wsComm.Range("C7").Font.Name = "Arial"
$wsComm.Range("C7").Font.Size = 10
$wsComm.Range("C7").Value = "ELC_GRID_RES"

$wsComm.Range("B8").Style = "Normal"
$wsComm.Range("B8").Value = "NRG"
$wsComm.Range("C8").Font.Name = "Arial"
$wsComm.Range("C8").Font.Size = 10
$wsComm.Range("C8").Value = "ELC_GRID"

$wsComm.Range("B9").Value = "DEM"
$wsComm.Range("C9").Value = "ELC_FIN"

$wsComm.Range("B10").Select()

# --- SEC_Processes sheet ---
$wsProc = $wb.Sheets.Item("SEC_Processes")

$wsProc.Range("B7").Style = "Normal"
$wsProc.Range("B7").Value = "DMD"
$wsProc.Range("D7").Font.Name = "Arial"
$wsProc.Range("D7").Font.Size = 10
$wsProc.Range("D7").Value = "ELC_FIN_DEM"

$wsProc.Activate()
$excel.ActiveWindow.Zoom = 110
$wsProc.Range("B8").Select()

# --- FINAL_DEMAD_PRC sheet ---
$wsFinal = $wb.Sheets.Item("FINAL_DEMAD_PRC")

$wsFinal.Range("B7").Formula = "=SEC_Processes!D7"
$wsFinal.Range("D8").Formula = "=SEC_Comm!C7"
$wsFinal.Range("D9").Formula = "=SEC_Comm!C8"
$wsFinal.Range("E10").Value = "ELC_FIN"

$wsFinal.Activate()
$wsFinal.Range("E11").Select()
